$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New customer rows are plain text (CPF/CEP/phone numbers must not turn into
# numbers), so force the incoming cells to Text format before writing values.
# Column D (Endereço) is skipped entirely for every new row, so it is left
# untouched here.
$ws.Range("A6:C15").NumberFormat = "@"
$ws.Range("E6:H15").NumberFormat = "@"

$data = @(
    @("roberto",        "123123124",    "991232",         "19123829",     "roberto@gmail.com",     "1928391823",     "Rua das tropaceiras"),
    @("bruno",           "12312",        "019230812938",   "819238192389", "brunofraga@gmail.com",  "182381723",      "ruas das molecas"),
    @("BRUNO DE FRAGA",  "12312",        "12312312",       "92320-195",    "qweqwe@qweqwe",         "132123123",      "Rua 3 Pinheiros I, 27"),
    @("joanues",         "1293i1923",    "-1203192309",    "123912839",    "joanues@gmail.com",     "123i912329",     "rua das horticias"),
    @("maria",           "1923192839",   "39139212068",    "91239182938",  "maria@gmail.com",       "1923891283",     "rua das alamedas"),
    @("ana luiza",       "SJAISJDI@@",   "19238912381273", "9123918239",   "analuiza@gmail.com",    "192u391239182",  "multi dimensoses"),
    @("leticia",         "98745451",     "91283918239",    "128391823981", "leticia@gmail.com",     "19283912839",    "rua das corticeiras"),
    @("joana",           "12381928391",  "1928391823",     "92320-195",    "joana@gmail.com",       "19823918293182", "Rua 3 Pinheiros I, 27"),
    @("francine",        "123912u3192",  "19283928391",    "192839123891", "fran@gmail.com",        "1923819283",     "rua tapajos"),
    @("joselito",        "12319238",     "1923891283",     "1923912839",   "qweqweqweq@qweqwe",     "91823918293",    "qjwdqwhduhq")
)

$row = 6
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    # column D (Endereço) is intentionally left blank for every new customer row
    $ws.Cells.Item($row, 5).Value = $rec[3]
    $ws.Cells.Item($row, 6).Value = $rec[4]
    $ws.Cells.Item($row, 7).Value = $rec[5]
    $ws.Cells.Item($row, 8).Value = $rec[6]
    $row++
}
